# Añadir link del servo
# Adds a "Link" column (G) to the "Inventario_cero_coste" sheet and a new
# row describing the SM-S2309S servo, including its datasheet/database link.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventario_cero_coste")

# --- New column G header ("Link", same text used in the other sheets) ---
$ws.Range("G1").Value = "Link"

# Give the new column a sensible width, matching the author's intent.
$ws.Columns.Item(7).ColumnWidth = 25.6

# --- New row (7) with the servo component info ---
$ws.Range("A7").Value = "SM-S2309S"
$ws.Range("C7").Value = "Servomotor, servo, motor"
$ws.Range("D7").Value = "Servo chiquito plasticucho de 150º de apertura. Tiene cable feedback."
$ws.Range("E7").Value = "Echedey"
$ws.Range("G7").Value = "https://servodatabase.com/servo/springrc/sm-s2309s"

# Reuse the header formatting (fill/style) for the new G1 cell, matching
# the look of the rest of the header row (F1).
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Link"

# Leave the selection on the newly added link cell, as in the final file.
$ws.Range("G7").Select()

$wb.Save()
